# Mifos Automation Excels / Client / 4014-CREATEACTIVECLIENT.xlsx
# "Share Product 15 Test Cases" - update the mobile number test value on the
# "Input" sheet and move the active selection to the cell that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# B6 ("mobilenumber" row) previously carried the input-field style (s="4")
# together with the old sample number. Reset it to the default "Normal"
# style and write the new sample mobile number.
$cell = $ws.Range("B6")
$cell.Style = "Normal"
$cell.Value = 9987654321

# Move/save the selection onto the cell that was just changed.
$cell.Select()
